$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $val) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value2 = $val
    $r.Style = "Normal"
}

Set-TextCell "B2" "Bitcoin"
Set-TextCell "C2" "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextCell "D2" "50.997.62"
Set-TextCell "E2" "  -0.86%  "

Set-TextCell "B3" "Ethereum"
Set-TextCell "C3" "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextCell "D3" "2.898.36"
Set-TextCell "E3" "  -0.71%  "

Set-TextCell "B4" "TetherUSD"
Set-TextCell "C4" "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextCell "D4" "1.00"
Set-TextCell "E4" "  +0.06%  "

Set-TextCell "B5" "BNB"
Set-TextCell "C5" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell "D5" "366.84"
Set-TextCell "E5" "  +4.89%  "

Set-TextCell "B6" "Solana"
Set-TextCell "C6" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell "D6" "101.88"
Set-TextCell "E6" "  -3.79%  "

Set-TextCell "B7" "XRP"
Set-TextCell "C7" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell "D7" "0.541"
Set-TextCell "E7" "  -2.30%  "

Set-TextCell "B8" "USDC"
Set-TextCell "C8" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell "D8" "1.00"
Set-TextCell "E8" "  +0.01%  "

Set-TextCell "B9" "Cardano"
Set-TextCell "C9" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell "D9" "0.579"
Set-TextCell "E9" "  -4.10%  "

Set-TextCell "B10" "Avalanche"
Set-TextCell "C10" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D10" "36.26"
Set-TextCell "E10" "  -3.80%  "

Set-TextCell "B11" "TRON"
Set-TextCell "C11" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D11" "0.139"
Set-TextCell "E11" "  +0.54%  "

Set-TextCell "B12" "Dogecoin"
Set-TextCell "C12" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell "D12" "0.0827"
Set-TextCell "E12" "  -2.26%  "

Set-TextCell "B13" "Chainlink"
Set-TextCell "C13" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D13" "18.12"
Set-TextCell "E13" "  -3.94%  "

Set-TextCell "B14" "WrappedliquidstakedEther2.0"
Set-TextCell "C14" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D14" "3.357.47"
Set-TextCell "E14" "  -0.43%  "

Set-TextCell "B15" "Polkadot"
Set-TextCell "C15" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D15" "7.32"
Set-TextCell "E15" "  -3.55%  "

Set-TextCell "B16" "WrappedEther"
Set-TextCell "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D16" "2.904.40"
Set-TextCell "E16" "  +0.26%  "

Set-TextCell "B17" "Polygon"
Set-TextCell "C17" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell "D17" "0.913"
Set-TextCell "E17" "  -4.60%  "

Set-TextCell "B18" "WrappedBTC"
Set-TextCell "C18" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D18" "50.964.69"
Set-TextCell "E18" "  -0.72%  "

Set-TextCell "B19" "ImmutableX"
Set-TextCell "C19" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D19" "3.20"
Set-TextCell "E19" "  -6.31%  "

Set-TextCell "B20" "Uniswap"
Set-TextCell "C20" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D20" "7.11"
Set-TextCell "E20" "  -3.93%  "

Set-TextCell "B21" "InternetComputer(DFINITY)"
Set-TextCell "C21" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D21" "12.77"
Set-TextCell "E21" "  -4.69%  "

Set-TextCell "B22" "ShibaInu"
Set-TextCell "C22" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D22" "0.0₃0937"
Set-TextCell "E22" "  -2.76%  "

Set-TextCell "B23" "Litecoin"
Set-TextCell "C23" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D23" "67.86"
Set-TextCell "E23" "  -1.41%  "

Set-TextCell "B24" "BitcoinCash"
Set-TextCell "C24" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D24" "257.49"
Set-TextCell "E24" "  -0.97%  "

Set-TextCell "B25" "PancakeSwap"
Set-TextCell "C25" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D25" "2.66"
Set-TextCell "E25" "  -1.14%  "

Set-TextCell "B26" "LEO"
Set-TextCell "C26" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D26" "4.33"
Set-TextCell "E26" "  +2.67%  "

Set-TextCell "B27" "Kaspa"
Set-TextCell "C27" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D27" "0.171"
Set-TextCell "E27" "  -1.26%  "

Set-TextCell "B28" "Dai"
Set-TextCell "C28" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D28" "1.00"
Set-TextCell "E28" "  +0.08%  "

Set-TextCell "B29" "EthereumClassic"
Set-TextCell "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D29" "25.45"
Set-TextCell "E29" "  -3.46%  "

Set-TextCell "B30" "Filecoin"
Set-TextCell "C30" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D30" "6.94"
Set-TextCell "E30" "  -5.90%  "

Set-TextCell "B31" "Hedera"
Set-TextCell "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D31" "0.101"
Set-TextCell "E31" "  -3.30%  "

Set-TextCell "B32" "RenderToken"
Set-TextCell "C32" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D32" "6.15"
Set-TextCell "E32" "  +1.58%  "

Set-TextCell "B33" "Cosmos"
Set-TextCell "C33" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D33" "9.83"
Set-TextCell "E33" "  -3.74%  "

Set-TextCell "B34" "Toncoin"
Set-TextCell "C34" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D34" "2.11"
Set-TextCell "E34" "  -1.34%  "

Set-TextCell "B35" "InjectiveProtocol"
Set-TextCell "C35" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D35" "34.14"
Set-TextCell "E35" "  -3.82%  "

Set-TextCell "B36" "OKB"
Set-TextCell "C36" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D36" "50.84"
Set-TextCell "E36" "  +1.00%  "

Set-TextCell "B37" "FirstDigitalUSD"
Set-TextCell "C37" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D37" "1.00"
Set-TextCell "E37" "  +0.42%  "

Set-TextCell "B38" "VeChain"
Set-TextCell "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D38" "0.0419"
Set-TextCell "E38" "  -1.46%  "

Set-TextCell "B39" "LidoDAOToken"
Set-TextCell "C39" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D39" "2.97"
Set-TextCell "E39" "  -5.00%  "

Set-TextCell "B40" "Stacks"
Set-TextCell "C40" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D40" "2.61"
Set-TextCell "E40" "  -1.13%  "

Set-TextCell "B41" "Celestia"
Set-TextCell "C41" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D41" "16.88"
Set-TextCell "E41" "  -4.09%  "

Set-TextCell "B42" "ARBITRUM"
Set-TextCell "C42" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D42" "1.82"
Set-TextCell "E42" "  -5.68%  "

Set-TextCell "B43" "Stellar"
Set-TextCell "C43" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D43" "0.112"
Set-TextCell "E43" "  -3.09%  "

Set-TextCell "B44" "EnergySwap"
Set-TextCell "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D44" "22.02"
Set-TextCell "E44" "  -1.54%  "

Set-TextCell "B45" "Monero"
Set-TextCell "C45" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D45" "118.34"
Set-TextCell "E45" "  -1.06%  "

Set-TextCell "B46" "WEMIXToken"
Set-TextCell "C46" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D46" "2.08"
Set-TextCell "E46" "  -1.99%  "

Set-TextCell "B47" "Maker"
Set-TextCell "C47" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D47" "2.011.57"
Set-TextCell "E47" "  -3.88%  "

Set-TextCell "B48" "ApeXProtocol"
Set-TextCell "C48" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell "D48" "2.31"
Set-TextCell "E48" "  +1.01%  "

Set-TextCell "B49" "NEARProtocol"
Set-TextCell "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D49" "3.13"
Set-TextCell "E49" "  -5.17%  "

Set-TextCell "B50" "RocketPoolETH"
Set-TextCell "C50" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell "D50" "3.186.54"
Set-TextCell "E50" "  -0.44%  "

Set-TextCell "B51" "TheGraph"
Set-TextCell "C51" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D51" "0.234"
Set-TextCell "E51" "  -1.72%  "
